$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted right after row 44, pushing every
# existing record from row 45 down one row (45->46, 46->47, ... 85->86).
$ws.Rows(45).Insert()

# Populate the newly-opened row 45 with the new record's data. The
# non-varying "template" columns (A,B,C,E,F,G,H,I,N,O,Q,R) are copied
# from the pattern shared by every other row in this table.
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = 44789
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100112040
$ws.Range("G45").Value = "Cilantro"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 300
$ws.Range("K45").Value = 2800
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = 2900
$ws.Range("N45").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 1450
$ws.Range("Q45").Value = 2
$ws.Range("R45").Value = "Hortaliza"
